$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 10
$ws.Range("B6").Value = "Juliana"
$ws.Range("C6").Value = "julichave09@gmail.com"
